$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: add "time_taken" column header with the same style as the
# other header cells (copy format from E1, which carries the bold/border/
# centered header style).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data cells F2:F73: per-row "time_taken" timestamps (no special style,
# matching the plain data cells in columns B-E).
$timestamps = @(
    "2021-10-05 13:42:02.101265",
    "2021-10-05 13:42:02.101278",
    "2021-10-05 13:42:02.101282",
    "2021-10-05 13:42:02.101285",
    "2021-10-05 13:42:02.101289",
    "2021-10-05 13:42:02.101292",
    "2021-10-05 13:42:02.101295",
    "2021-10-05 13:42:02.101298",
    "2021-10-05 13:42:02.101301",
    "2021-10-05 13:42:02.101304",
    "2021-10-05 13:42:02.101307",
    "2021-10-05 13:42:02.101310",
    "2021-10-05 13:42:02.101313",
    "2021-10-05 13:42:02.101316",
    "2021-10-05 13:42:02.101319",
    "2021-10-05 13:42:02.101322",
    "2021-10-05 13:42:02.101326",
    "2021-10-05 13:42:02.101329",
    "2021-10-05 13:42:02.101332",
    "2021-10-05 13:42:02.101335",
    "2021-10-05 13:42:02.101338",
    "2021-10-05 13:42:02.101341",
    "2021-10-05 13:42:02.101344",
    "2021-10-05 13:42:02.101347",
    "2021-10-05 13:42:02.101350",
    "2021-10-05 13:42:02.101353",
    "2021-10-05 13:42:02.101356",
    "2021-10-05 13:42:02.101359",
    "2021-10-05 13:42:02.101362",
    "2021-10-05 13:42:02.101365",
    "2021-10-05 13:42:02.101368",
    "2021-10-05 13:42:02.101371",
    "2021-10-05 13:42:02.101375",
    "2021-10-05 13:42:02.101378",
    "2021-10-05 13:42:02.101381",
    "2021-10-05 13:42:02.101384",
    "2021-10-05 13:42:02.101387",
    "2021-10-05 13:42:02.101389",
    "2021-10-05 13:42:02.101392",
    "2021-10-05 13:42:02.101395",
    "2021-10-05 13:42:02.101399",
    "2021-10-05 13:42:02.101402",
    "2021-10-05 13:42:02.101405",
    "2021-10-05 13:42:02.101408",
    "2021-10-05 13:42:02.101411",
    "2021-10-05 13:42:02.101414",
    "2021-10-05 13:42:02.101417",
    "2021-10-05 13:42:02.101420",
    "2021-10-05 13:42:02.101423",
    "2021-10-05 13:42:02.101425",
    "2021-10-05 13:42:02.101428",
    "2021-10-05 13:42:02.101431",
    "2021-10-05 13:42:02.101435",
    "2021-10-05 13:42:02.101438",
    "2021-10-05 13:42:02.101441",
    "2021-10-05 13:42:02.101444",
    "2021-10-05 13:42:02.101447",
    "2021-10-05 13:42:02.101450",
    "2021-10-05 13:42:02.101452",
    "2021-10-05 13:42:02.101456",
    "2021-10-05 13:42:02.101459",
    "2021-10-05 13:42:02.101462",
    "2021-10-05 13:42:02.101465",
    "2021-10-05 13:42:02.101468",
    "2021-10-05 13:42:02.101472",
    "2021-10-05 13:42:02.101476",
    "2021-10-05 13:42:02.101478",
    "2021-10-05 13:42:02.101481",
    "2021-10-05 13:42:02.101484",
    "2021-10-05 13:42:02.101487",
    "2021-10-05 13:42:02.101490",
    "2021-10-05 13:42:02.101493"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
